$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cells per diff
$ws.Range("B2").Value = "Open Source"
$ws.Range("B3").Value = "Under Review"

$ws.Range("A5").Value = "/src/api/http-client.json"
$ws.Range("B5").Value = "Pending Review"
$ws.Range("C5").Value = "Axios HTTP Client"

$ws.Range("A6").Value = "/src/utils/date-formatter.js"
$ws.Range("B6").Value = "Open Source"
$ws.Range("C6").Value = "Moment.js"

# Add new row 7
$ws.Range("A7").Value = "/src/components/charts/bar-chart.css"
$ws.Range("B7").Value = "Approved"
$ws.Range("C7").Value = "Chart.js"
